# edit.ps1 - Updates the NYPD 113th Precinct weekly CompStat report
# Commit: "New crime data collected"
#
# This script mutates the open workbook ($excel.ActiveWorkbook) to reflect
# the newly collected week's crime statistics: the report header (volume
# number + reporting week dates) and the full crime-complaint stats table
# (rows 15-33, columns C:N -- Week to Date / 28 Day / Year to Date / 2 Year
# / 15 Year / 32 Year figures and their percent changes).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates -------------------------------------------------
# "Volume 32   Number  5" -> "Volume 32   Number  6"
$ws.Range("A8").Value = "Volume 32   Number  6"

# "Report Covering the Week  1/27/2025  Through  2/2/2025"
#   -> "Report Covering the Week  2/3/2025  Through  2/9/2025"
$ws.Range("C9").Value = "Report Covering the Week  2/3/2025  Through  2/9/2025"

# --- simple numeric updates (no format/type change) ---
$simpleValues = @{
    "G15" = 1
    "H15" = 100
    "N15" = -84.615384615384
    "C16" = 1
    "F16" = 5
    "G16" = 11
    "H16" = -54.545454545454
    "I16" = 10
    "J16" = 15
    "K16" = -33.333333333333
    "L16" = 0
    "M16" = -72.222222222222
    "N16" = -91.525423728813
    "C17" = 9
    "D17" = 5
    "E17" = 80
    "F17" = 31
    "G17" = 25
    "H17" = 24
    "I17" = 40
    "J17" = 38
    "K17" = 5.263157894736
    "L17" = 2.564102564102
    "N17" = -45.205479452054
    "D18" = 1
    "E18" = 0
    "F18" = 4
    "H18" = -20
    "I18" = 5
    "J18" = 6
    "K18" = -16.666666666666
    "L18" = -16.666666666666
    "M18" = -88.888888888888
    "N18" = -94.791666666666
    "C19" = 9
    "D19" = 8
    "E19" = 12.5
    "G19" = 22
    "H19" = 40.909090909090
    "I19" = 42
    "J19" = 29
    "K19" = 44.827586206896
    "L19" = 31.25
    "M19" = -39.130434782608
    "N19" = -78.238341968911
    "C20" = 3
    "D20" = 7
    "E20" = -57.142857142857
    "F20" = 11
    "G20" = 14
    "H20" = -21.428571428571
    "I20" = 22
    "J20" = 24
    "K20" = -8.333333333333
    "L20" = 69.230769230769
    "M20" = -31.25
    "N20" = -88.541666666666
    "C21" = 23
    "D21" = 26
    "E21" = -11.538461538461
    "F21" = 85
    "G21" = 78
    "H21" = 8.974358974358
    "I21" = 122
    "J21" = 114
    "K21" = 7.017543859649
    "L21" = 19.607843137254
    "M21" = -42.452830188679
    "N21" = -82.293178519593
    "G23" = 1
    "D24" = 26
    "E24" = -46.153846153846
    "F24" = 68
    "G24" = 73
    "H24" = -6.849315068493
    "I24" = 87
    "J24" = 115
    "K24" = -24.347826086956
    "L24" = -23.684210526315
    "M24" = -9.375
    "C25" = 5
    "D25" = 7
    "E25" = -28.571428571428
    "F25" = 17
    "G25" = 19
    "H25" = -10.526315789473
    "I25" = 22
    "J25" = 42
    "K25" = -47.619047619047
    "L25" = -31.25
    "C26" = 8
    "D26" = 16
    "E26" = -50
    "F26" = 36
    "G26" = 44
    "H26" = -18.181818181818
    "I26" = 61
    "J26" = 60
    "K26" = 1.666666666666
    "L26" = 12.962962962963
    "M26" = -15.277777777777
    "E27" = -100
    "J27" = 3
    "K27" = -33.333333333333
    "C28" = 1
    "G28" = 2
    "H28" = 300
    "I28" = 9
    "K28" = 80
    "L28" = -18.181818181818
    "L29" = 0
    "M29" = -25
    "L30" = 50
    "M30" = 0
    "L33" = 0
}
foreach ($ref in $simpleValues.Keys) {
    $ws.Range($ref).Value = $simpleValues[$ref]
}

# --- cells converting to text placeholders ('0' / '***.*') ---
$textValues = @{
    "C15" = "0"
    "D15" = "0"
    "E15" = "***.*"
    "D23" = "0"
    "E23" = "***.*"
    "C27" = "0"
    "C29" = "0"
    "D29" = "0"
    "E29" = "***.*"
    "C30" = "0"
    "D30" = "0"
    "E30" = "***.*"
}
foreach ($ref in $textValues.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $textValues[$ref]
}

# --- cells converting from text placeholders back to numbers ---
$numberFixups = @(
    @{ Ref = "D16"; Value = 5; Format = '#,##0' }
    @{ Ref = "E16"; Value = -80; Format = '#,##0.0;"-"#,##0.0' }
    @{ Ref = "C18"; Value = 1; Format = '#,##0' }
    @{ Ref = "F33"; Value = 1; Format = '#,##0' }
    @{ Ref = "I33"; Value = 1; Format = '#,##0' }
)
foreach ($item in $numberFixups) {
    $cell = $ws.Range($item.Ref)
    $cell.NumberFormat = $item.Format
    $cell.Value = $item.Value
}
